$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows.Item(10).Delete()
